$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new header "Save" in H1, copying formatting (style) from the
# neighboring header cell G1 so it matches the other header cells.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill H2:H5 with the value 1
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 1
